# Updated cryptos list with latest prices / 1h volume changes.
# Column D holds the "Price" text and Column E the "Volume(1h)" text for
# each coin row. A handful of rows also show the coin list being
# re-ranked (MultiversX and LidoDAOToken swap places, rows 40/41).
#
# Because several price strings are plain numeric text (e.g. "243.80",
# "0.628", ...), Excel's normal type-inference would silently convert
# them into numbers (dropping trailing zeros / changing representation)
# if assigned directly. To keep them as text - matching the source data
# - those values are entered with a leading apostrophe (Excel's
# "force text" marker) and the cell style is then reset back to Normal
# so no stray text-number-format style lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param($addr, $text)
    # Plain decimals (e.g. "243.80", "0.0962") would otherwise be
    # auto-converted to numbers by Excel's type inference; values with
    # more than one "." (e.g. "42.277.59") already round-trip as text.
    if ($text -match '^[0-9]+(\.[0-9]+)?$') {
        $ws.Range($addr).Value = "'" + $text
        $ws.Range($addr).Style = "Normal"
    } else {
        $ws.Range($addr).Value = $text
    }
}

# --- Row 2: Bitcoin ---
Set-PriceText "D2" "42.277.59"
$ws.Range("E2").Value = "  -0.81%  "

# --- Row 3: Ethereum ---
Set-PriceText "D3" "2.236.59"
$ws.Range("E3").Value = "  -0.65%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  -0.13%  "

# --- Row 5: BNB ---
Set-PriceText "D5" "243.80"
$ws.Range("E5").Value = "  -0.98%  "

# --- Row 6: XRP ---
Set-PriceText "D6" "0.628"
$ws.Range("E6").Value = "  -0.57%  "

# --- Row 7 ---
Set-PriceText "D7" "74.69"
$ws.Range("E7").Value = "  -2.51%  "

# --- Row 8 ---
$ws.Range("E8").Value = "  -0.08%  "

# --- Row 9 ---
Set-PriceText "D9" "0.606"
$ws.Range("E9").Value = "  -3.63%  "

# --- Row 10 ---
Set-PriceText "D10" "42.78"
$ws.Range("E10").Value = "  -5.67%  "

# --- Row 11 ---
Set-PriceText "D11" "0.0962"
$ws.Range("E11").Value = "  +0.97%  "

# --- Row 12 ---
Set-PriceText "D12" "7.02"
$ws.Range("E12").Value = "  -4.23%  "

# --- Row 13 ---
Set-PriceText "D13" "0.104"
$ws.Range("E13").Value = "  +1.03%  "

# --- Row 14 ---
Set-PriceText "D14" "2.565.84"
$ws.Range("E14").Value = "  -1.58%  "

# --- Row 15 ---
Set-PriceText "D15" "14.38"
$ws.Range("E15").Value = "  -2.37%  "

# --- Row 16 ---
Set-PriceText "D16" "0.842"
$ws.Range("E16").Value = "  -2.61%  "

# --- Row 17 ---
Set-PriceText "D17" "2.226.13"
$ws.Range("E17").Value = "  -1.40%  "

# --- Row 18 ---
Set-PriceText "D18" "42.036.95"
$ws.Range("E18").Value = "  -0.98%  "

# --- Row 19 ---
Set-PriceText "D19" "0.0000107"
$ws.Range("E19").Value = "  +5.09%  "

# --- Row 20 ---
Set-PriceText "D20" "6.25"
$ws.Range("E20").Value = "  +0.50%  "

# --- Row 21 ---
Set-PriceText "D21" "73.08"
$ws.Range("E21").Value = "  +1.23%  "

# --- Row 22 ---
Set-PriceText "D22" "11.29"
$ws.Range("E22").Value = "  +0.80%  "

# --- Row 23 ---
Set-PriceText "D23" "231.31"
$ws.Range("E23").Value = "  -0.46%  "

# --- Row 24 ---
Set-PriceText "D24" "2.10"
$ws.Range("E24").Value = "  -7.10%  "

# --- Row 25 ---
$ws.Range("E25").Value = "  +0.12%  "

# --- Row 26 ---
Set-PriceText "D26" "11.46"
$ws.Range("E26").Value = "  -3.93%  "

# --- Row 27 ---
$ws.Range("E27").Value = "  -0.06%  "

# --- Row 28 ---
Set-PriceText "D28" "2.27"
$ws.Range("E28").Value = "  -2.05%  "

# --- Row 29 ---
$ws.Range("E29").Value = "  -3.56%  "

# --- Row 30 ---
Set-PriceText "D30" "167.07"
$ws.Range("E30").Value = "  -0.18%  "

# --- Row 31 ---
Set-PriceText "D31" "20.63"
$ws.Range("E31").Value = "  -0.36%  "

# --- Row 32 ---
Set-PriceText "D32" "5.68"
$ws.Range("E32").Value = "  +0.99%  "

# --- Row 33 ---
Set-PriceText "D33" "0.0803"
$ws.Range("E33").Value = "  -3.12%  "

# --- Row 34 ---
Set-PriceText "D34" "30.04"
$ws.Range("E34").Value = "  -6.77%  "

# --- Row 35 ---
$ws.Range("E35").Value = "  -0.39%  "

# --- Row 36 ---
Set-PriceText "D36" "0.110"
$ws.Range("E36").Value = "  -8.63%  "

# --- Row 37 ---
Set-PriceText "D37" "4.35"
$ws.Range("E37").Value = "  -6.96%  "

# --- Row 38 ---
$ws.Range("E38").Value = "  -3.89%  "

# --- Row 39 ---
Set-PriceText "D39" "13.49"
$ws.Range("E39").Value = "  -6.73%  "

# --- Rows 40 & 41: MultiversX and LidoDAOToken swap ranking order ---
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-PriceText "D40" "2.14"
$ws.Range("E40").Value = "  -2.41%  "

$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-PriceText "D41" "65.49"
$ws.Range("E41").Value = "  +1.92%  "

# --- Row 42 ---
Set-PriceText "D42" "5.71"
$ws.Range("E42").Value = "  -1.98%  "

# --- Row 43 ---
Set-PriceText "D43" "0.199"
$ws.Range("E43").Value = "  -2.23%  "

# --- Row 44 ---
$ws.Range("E44").Value = "  -2.00%  "

# --- Row 45 ---
Set-PriceText "D45" "104.88"
$ws.Range("E45").Value = "  -3.50%  "

# --- Row 46 ---
Set-PriceText "D46" "0.100"
$ws.Range("E46").Value = "  -3.01%  "

# --- Row 47 ---
Set-PriceText "D47" "2.37"
$ws.Range("E47").Value = "  -2.27%  "

# --- Row 48 ---
$ws.Range("E48").Value = "  -1.84%  "

# --- Row 49 ---
Set-PriceText "D49" "1.18"
$ws.Range("E49").Value = "  -1.34%  "

# --- Row 50 ---
$ws.Range("E50").Value = "  -1.12%  "

# --- Row 51 ---
Set-PriceText "D51" "2.439.46"
$ws.Range("E51").Value = "  -1.29%  "
